$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-13: column B -> -1, column C -> "date"
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 2).Value = -1
    $ws.Cells.Item($r, 3).Value = "date"
}

# Specific U-column (csim) updates
$ws.Range("U3").Value = 0
$ws.Range("U8").Value = 0
$ws.Range("U9").Value = 0

# New row 14 - metadata for an additional attribute discovered in the dataset
$ws.Range("A14").Value = ""
$ws.Range("B14").Value = -1
$ws.Range("C14").Value = "date"
$ws.Range("D14").Value = 177471
$ws.Range("E14").Value = 177471
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = -1
$ws.Range("L14").Value = -1
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 9
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 177455
$ws.Range("Q14").Value = 0
$ws.Range("R14").Value = "<Unspecified>"
$ws.Range("S14").Value = ""
$ws.Range("T14").Value = 177455
$ws.Range("U14").Value = 0
$ws.Range("V14").Value = 113
$ws.Range("W14").Value = 4
$ws.Range("X14").Value = "['GILLESPIE', 'IZETT', 'PATERSON', 'PITCHER']"
